# Update Pais (COVID-19 countries) worksheet with refreshed case figures
# and the two country-name swaps that the new sort order implies
# (rows keep themselves sorted by "Casos totales" descending).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: Datos actualizados a 19 de Mayo de 2020 a las 21:05
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 19 de Mayo de 2020 a las 21:05'

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 1561399
$ws.Cells.Item(4, 3).Value = 11105
$ws.Cells.Item(4, 4).Value = 360088
$ws.Cells.Item(4, 5).Value = 1108507
$ws.Cells.Item(4, 7).Value = 823
$ws.Cells.Item(4, 8).Value = 92804

# Row 6: España
$ws.Cells.Item(6, 2).Value = 278803
$ws.Cells.Item(6, 3).Value = 615
$ws.Cells.Item(6, 5).Value = 54067
$ws.Cells.Item(6, 7).Value = 69
$ws.Cells.Item(6, 8).Value = 27778

# Row 10: Francia
$ws.Cells.Item(10, 2).Value = 180809
$ws.Cells.Item(10, 3).Value = 882
$ws.Cells.Item(10, 4).Value = 62563
$ws.Cells.Item(10, 5).Value = 90007

# Row 11: Alemania
$ws.Cells.Item(11, 2).Value = 177739
$ws.Cells.Item(11, 3).Value = 450
$ws.Cells.Item(11, 5).Value = 13870
$ws.Cells.Item(11, 7).Value = 46
$ws.Cells.Item(11, 8).Value = 8169

# Row 14: India
$ws.Cells.Item(14, 2).Value = 106468
$ws.Cells.Item(14, 3).Value = 6140
$ws.Cells.Item(14, 5).Value = 60860

# Row 52: Argentina
$ws.Cells.Item(52, 4).Value = 2872
$ws.Cells.Item(52, 5).Value = 5117

# Row 90: Consejo Danes para los Refugiados
$ws.Cells.Item(90, 1).Value = 'Consejo Danes para los Refugiados'
$ws.Cells.Item(90, 2).Value = 1629
$ws.Cells.Item(90, 3).Value = 91
$ws.Cells.Item(90, 4).Value = 290
$ws.Cells.Item(90, 5).Value = 1278
$ws.Cells.Item(90, 7).Value = 0
$ws.Cells.Item(90, 8).Value = 61

# Row 91: Lituania
$ws.Cells.Item(91, 1).Value = 'Lituania'
$ws.Cells.Item(91, 2).Value = 1562
$ws.Cells.Item(91, 3).Value = 15
$ws.Cells.Item(91, 4).Value = 1025
$ws.Cells.Item(91, 5).Value = 477
$ws.Cells.Item(91, 7).Value = 1
$ws.Cells.Item(91, 8).Value = 60

# Row 123: Jordania
$ws.Cells.Item(123, 2).Value = 649
$ws.Cells.Item(123, 3).Value = 20
$ws.Cells.Item(123, 5).Value = 227

# Row 134: Nepal
$ws.Cells.Item(134, 2).Value = 407
$ws.Cells.Item(134, 3).Value = 32
$ws.Cells.Item(134, 5).Value = 368

# Row 161: Mauritania
$ws.Cells.Item(161, 1).Value = 'Mauritania'
$ws.Cells.Item(161, 2).Value = 131
$ws.Cells.Item(161, 3).Value = 50
$ws.Cells.Item(161, 4).Value = 7
$ws.Cells.Item(161, 5).Value = 120
$ws.Cells.Item(161, 8).Value = 4

# Row 162: Yemen
$ws.Cells.Item(162, 1).Value = 'Yemen'
$ws.Cells.Item(162, 2).Value = 130
$ws.Cells.Item(162, 4).Value = 1
$ws.Cells.Item(162, 5).Value = 109
$ws.Cells.Item(162, 8).Value = 20

# Row 163: Bermudas
$ws.Cells.Item(163, 1).Value = 'Bermudas'
$ws.Cells.Item(163, 2).Value = 125
$ws.Cells.Item(163, 4).Value = 77
$ws.Cells.Item(163, 5).Value = 39
$ws.Cells.Item(163, 8).Value = 9

# Row 164: Guyana
$ws.Cells.Item(164, 1).Value = 'Guyana'
$ws.Cells.Item(164, 2).Value = 124
$ws.Cells.Item(164, 4).Value = 45
$ws.Cells.Item(164, 5).Value = 69
$ws.Cells.Item(164, 8).Value = 10

# Row 165: Camboya
$ws.Cells.Item(165, 1).Value = 'Camboya'
$ws.Cells.Item(165, 2).Value = 122
$ws.Cells.Item(165, 4).Value = 122
$ws.Cells.Item(165, 5).Value = 0
$ws.Cells.Item(165, 8).Value = 0

# Row 166: Trinidad yTobago
$ws.Cells.Item(166, 1).Value = 'Trinidad yTobago'
$ws.Cells.Item(166, 2).Value = 116
$ws.Cells.Item(166, 4).Value = 107
$ws.Cells.Item(166, 5).Value = 1
$ws.Cells.Item(166, 8).Value = 8

# Row 167: Aruba
$ws.Cells.Item(167, 1).Value = 'Aruba'
$ws.Cells.Item(167, 2).Value = 101
$ws.Cells.Item(167, 4).Value = 93
$ws.Cells.Item(167, 5).Value = 5
$ws.Cells.Item(167, 8).Value = 3

# Row 168: Monaco
$ws.Cells.Item(168, 1).Value = 'Monaco'
$ws.Cells.Item(168, 2).Value = 97
$ws.Cells.Item(168, 4).Value = 87
$ws.Cells.Item(168, 5).Value = 6
$ws.Cells.Item(168, 8).Value = 4

# Row 169: Bahamas
$ws.Cells.Item(169, 1).Value = 'Bahamas'
$ws.Cells.Item(169, 2).Value = 96
$ws.Cells.Item(169, 4).Value = 43
$ws.Cells.Item(169, 5).Value = 42
$ws.Cells.Item(169, 8).Value = 11

# Row 170: Islas Caimanes
$ws.Cells.Item(170, 1).Value = 'Islas Caimanes'
$ws.Cells.Item(170, 2).Value = 94
$ws.Cells.Item(170, 4).Value = 55
$ws.Cells.Item(170, 5).Value = 38
$ws.Cells.Item(170, 8).Value = 1

# Row 171: Barbados
$ws.Cells.Item(171, 1).Value = 'Barbados'
$ws.Cells.Item(171, 2).Value = 88
$ws.Cells.Item(171, 4).Value = 68
$ws.Cells.Item(171, 5).Value = 13
$ws.Cells.Item(171, 8).Value = 7

# Row 172: Liechtenstein
$ws.Cells.Item(172, 1).Value = 'Liechtenstein'
$ws.Cells.Item(172, 2).Value = 82
$ws.Cells.Item(172, 4).Value = 55
$ws.Cells.Item(172, 5).Value = 26
$ws.Cells.Item(172, 8).Value = 1

# Row 196: Nueva Caledonia
$ws.Cells.Item(196, 1).Value = 'Nueva Caledonia'
$ws.Cells.Item(196, 4).Value = 18
$ws.Cells.Item(196, 8).Value = 0

# Row 197: Belice
$ws.Cells.Item(197, 1).Value = 'Belice'
$ws.Cells.Item(197, 4).Value = 16
$ws.Cells.Item(197, 8).Value = 2

# Row 209: Seychelles
$ws.Cells.Item(209, 1).Value = 'Seychelles'

# Row 210: Groenlandia
$ws.Cells.Item(210, 1).Value = 'Groenlandia'
$ws.Cells.Item(210, 4).Value = 11
$ws.Cells.Item(210, 8).Value = 0

# Row 211: Montserrat
$ws.Cells.Item(211, 1).Value = 'Montserrat'
$ws.Cells.Item(211, 4).Value = 10
$ws.Cells.Item(211, 8).Value = 1

# Row 215: San Bartolome
$ws.Cells.Item(215, 1).Value = 'San Bartolome'

# Row 216: Bonaire, San Eustaquio y Saba
$ws.Cells.Item(216, 1).Value = 'Bonaire, San Eustaquio y Saba'
